$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data blocks between rows 2-5 (originally date 44908) and
# rows 6-9 (originally date 44890), matching row-by-row: 2<->6, 3<->7,
# 4<->8, 5<->9. Columns involved: D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen), S (Precio $/Kg).

$pairs = @(
    @{ a = 2; b = 6 },
    @{ a = 3; b = 7 },
    @{ a = 4; b = 8 },
    @{ a = 5; b = 9 }
)

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($pair in $pairs) {
    $rowA = $pair.a
    $rowB = $pair.b

    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")

        $valA = $rangeA.Value2
        $valB = $rangeB.Value2

        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}
